$wb = $excel.ActiveWorkbook

# ================= Sheet 1: Overall Scores =================
$ws = $wb.Worksheets.Item(1)
$blank = $ws.Cells.Item(500, 500)   # untouched cell, used as a blank Copy() source

$ws.Cells.Item(1, 1).Value = "Image Name"
$ws.Cells.Item(1, 2).Value = "Full Image Gen"
$ws.Cells.Item(1, 3).Value = "Segment Gen"

$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 40
$ws.Cells.Item(2, 3).Value = 80

$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 67
$ws.Cells.Item(3, 3).Value = 75

$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 70
$ws.Cells.Item(4, 3).Value = 70

$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 80
$ws.Cells.Item(5, 3).Value = 20

$ws.Cells.Item(6, 1).Value = "ecommerce"
$ws.Cells.Item(6, 2).Value = 95
$ws.Cells.Item(6, 3).Value = 93

$ws.Cells.Item(7, 1).Value = "webflow-full"
$ws.Cells.Item(7, 2).Value = 60
$ws.Cells.Item(7, 3).Value = 50

$ws.Cells.Item(8, 1).Value = "complex"
$ws.Cells.Item(8, 2).Value = 90
$ws.Cells.Item(8, 3).Value = 90

$ws.Cells.Item(9, 1).Value = "Superlist_website"
$ws.Cells.Item(9, 2).Value = 83
$ws.Cells.Item(9, 3).Value = 82

$ws.Cells.Item(10, 1).Value = "Spotify_website"
$ws.Cells.Item(10, 2).Value = 77
$ws.Cells.Item(10, 3).Value = 60

$ws.Cells.Item(11, 1).Value = "blog"
$ws.Cells.Item(11, 2).Value = 93
$ws.Cells.Item(11, 3).Value = 93

$ws.Cells.Item(12, 1).Value = "mubasic_website"
$ws.Cells.Item(12, 2).Value = 95
$ws.Cells.Item(12, 3).Value = 87

$ws.Cells.Item(13, 1).Value = "overflow_website"
$ws.Cells.Item(13, 2).Value = 87
$ws.Cells.Item(13, 3).Value = 83

$ws.Cells.Item(14, 1).Value = "RCA_website"
$ws.Cells.Item(14, 2).Value = 30
$ws.Cells.Item(14, 3).Value = 60

$ws.Cells.Item(15, 1).Value = "Crypto_website"
$ws.Cells.Item(15, 2).Value = 87
$ws.Cells.Item(15, 3).Value = 95

$ws.Cells.Item(16, 2).Value = 87
$ws.Cells.Item(16, 3).Value = 92
$ws.Cells.Item(16, 4).Value = "Crypto_website"

# Blank placeholder cells: exist in the sheet but carry no value, mirroring
# the authored file's empty inlineStr cells.
for ($r = 2; $r -le 15; $r++) {
    $blank.Copy($ws.Cells.Item($r, 4))
}
$blank.Copy($ws.Cells.Item(16, 1))

# D1 is blank too, but keeps the bold/bordered header style ("s=1") - copy it
# from A1 (same style) before A1 had a value, i.e. restore A1 after copying.
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 4))
$ws.Cells.Item(1, 4).Value = ""

# ================= Sheet 2: Visual Scores =================
$ws = $wb.Worksheets.Item(2)
$blank = $ws.Cells.Item(500, 500)   # untouched cell, used as a blank Copy() source

$ws.Cells.Item(1, 1).Value = "Image Name"
$ws.Cells.Item(1, 2).Value = "Full Image Gen"
$ws.Cells.Item(1, 3).Value = "Segment Gen"

$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 50
$ws.Cells.Item(2, 3).Value = 70

$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 50
$ws.Cells.Item(3, 3).Value = 60

$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 60
$ws.Cells.Item(4, 3).Value = 60

$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 80
$ws.Cells.Item(5, 3).Value = 30

$ws.Cells.Item(6, 1).Value = "ecommerce"
$ws.Cells.Item(6, 2).Value = 95
$ws.Cells.Item(6, 3).Value = 95

$ws.Cells.Item(7, 1).Value = "webflow-full"
$ws.Cells.Item(7, 2).Value = 60
$ws.Cells.Item(7, 3).Value = 60

$ws.Cells.Item(8, 1).Value = "complex"
$ws.Cells.Item(8, 2).Value = 80
$ws.Cells.Item(8, 3).Value = 85

$ws.Cells.Item(9, 1).Value = "Superlist_website"
$ws.Cells.Item(9, 2).Value = 85
$ws.Cells.Item(9, 3).Value = 85

$ws.Cells.Item(10, 1).Value = "Spotify_website"
$ws.Cells.Item(10, 2).Value = 75
$ws.Cells.Item(10, 3).Value = 70

$ws.Cells.Item(11, 1).Value = "blog"
$ws.Cells.Item(11, 2).Value = 95
$ws.Cells.Item(11, 3).Value = 95

$ws.Cells.Item(12, 1).Value = "mubasic_website"
$ws.Cells.Item(12, 2).Value = 95
$ws.Cells.Item(12, 3).Value = 80

$ws.Cells.Item(13, 1).Value = "overflow_website"
$ws.Cells.Item(13, 2).Value = 85
$ws.Cells.Item(13, 3).Value = 85

$ws.Cells.Item(14, 1).Value = "RCA_website"
$ws.Cells.Item(14, 2).Value = 40
$ws.Cells.Item(14, 3).Value = 60

$ws.Cells.Item(15, 1).Value = "Crypto_website"
$ws.Cells.Item(15, 2).Value = 85
$ws.Cells.Item(15, 3).Value = 95

$ws.Cells.Item(16, 2).Value = 85
$ws.Cells.Item(16, 3).Value = 95
$ws.Cells.Item(16, 4).Value = "Crypto_website"

# Blank placeholder cells: exist in the sheet but carry no value, mirroring
# the authored file's empty inlineStr cells.
for ($r = 2; $r -le 15; $r++) {
    $blank.Copy($ws.Cells.Item($r, 4))
}
$blank.Copy($ws.Cells.Item(16, 1))

# D1 is blank too, but keeps the bold/bordered header style ("s=1") - copy it
# from A1 (same style) before A1 had a value, i.e. restore A1 after copying.
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 4))
$ws.Cells.Item(1, 4).Value = ""

# ================= Sheet 3: Content Scores =================
$ws = $wb.Worksheets.Item(3)
$blank = $ws.Cells.Item(500, 500)   # untouched cell, used as a blank Copy() source

$ws.Cells.Item(1, 1).Value = "Image Name"
$ws.Cells.Item(1, 2).Value = "Full Image Gen"
$ws.Cells.Item(1, 3).Value = "Segment Gen"

$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 40
$ws.Cells.Item(2, 3).Value = 90

$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 90
$ws.Cells.Item(3, 3).Value = 90

$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 80
$ws.Cells.Item(4, 3).Value = 85

$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 90
$ws.Cells.Item(5, 3).Value = 20

$ws.Cells.Item(6, 1).Value = "ecommerce"
$ws.Cells.Item(6, 2).Value = 100
$ws.Cells.Item(6, 3).Value = 95

$ws.Cells.Item(7, 1).Value = "webflow-full"
$ws.Cells.Item(7, 2).Value = 70
$ws.Cells.Item(7, 3).Value = 50

$ws.Cells.Item(8, 1).Value = "complex"
$ws.Cells.Item(8, 2).Value = 100
$ws.Cells.Item(8, 3).Value = 95

$ws.Cells.Item(9, 1).Value = "Superlist_website"
$ws.Cells.Item(9, 2).Value = 90
$ws.Cells.Item(9, 3).Value = 90

$ws.Cells.Item(10, 1).Value = "Spotify_website"
$ws.Cells.Item(10, 2).Value = 85
$ws.Cells.Item(10, 3).Value = 60

$ws.Cells.Item(11, 1).Value = "blog"
$ws.Cells.Item(11, 2).Value = 95
$ws.Cells.Item(11, 3).Value = 95

$ws.Cells.Item(12, 1).Value = "mubasic_website"
$ws.Cells.Item(12, 2).Value = 100
$ws.Cells.Item(12, 3).Value = 95

$ws.Cells.Item(13, 1).Value = "overflow_website"
$ws.Cells.Item(13, 2).Value = 95
$ws.Cells.Item(13, 3).Value = 90

$ws.Cells.Item(14, 1).Value = "RCA_website"
$ws.Cells.Item(14, 2).Value = 30
$ws.Cells.Item(14, 3).Value = 70

$ws.Cells.Item(15, 1).Value = "Crypto_website"
$ws.Cells.Item(15, 2).Value = 95
$ws.Cells.Item(15, 3).Value = 100

$ws.Cells.Item(16, 2).Value = 95
$ws.Cells.Item(16, 3).Value = 95
$ws.Cells.Item(16, 4).Value = "Crypto_website"

# Blank placeholder cells: exist in the sheet but carry no value, mirroring
# the authored file's empty inlineStr cells.
for ($r = 2; $r -le 15; $r++) {
    $blank.Copy($ws.Cells.Item($r, 4))
}
$blank.Copy($ws.Cells.Item(16, 1))

# D1 is blank too, but keeps the bold/bordered header style ("s=1") - copy it
# from A1 (same style) before A1 had a value, i.e. restore A1 after copying.
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 4))
$ws.Cells.Item(1, 4).Value = ""

# ================= Sheet 4: Functional Scores =================
$ws = $wb.Worksheets.Item(4)
$blank = $ws.Cells.Item(500, 500)   # untouched cell, used as a blank Copy() source

$ws.Cells.Item(1, 1).Value = "Image Name"
$ws.Cells.Item(1, 2).Value = "Full Image Gen"
$ws.Cells.Item(1, 3).Value = "Segment Gen"

$ws.Cells.Item(2, 1).Value = "Figma_food_website"
$ws.Cells.Item(2, 2).Value = 30
$ws.Cells.Item(2, 3).Value = 80

$ws.Cells.Item(3, 1).Value = "Trees_website"
$ws.Cells.Item(3, 2).Value = 60
$ws.Cells.Item(3, 3).Value = 70

$ws.Cells.Item(4, 1).Value = "JoeCoffee-website-february-2022"
$ws.Cells.Item(4, 2).Value = 70
$ws.Cells.Item(4, 3).Value = 75

$ws.Cells.Item(5, 1).Value = "Hyer_website"
$ws.Cells.Item(5, 2).Value = 70
$ws.Cells.Item(5, 3).Value = 10

$ws.Cells.Item(6, 1).Value = "ecommerce"
$ws.Cells.Item(6, 2).Value = 90
$ws.Cells.Item(6, 3).Value = 90

$ws.Cells.Item(7, 1).Value = "webflow-full"
$ws.Cells.Item(7, 2).Value = 50
$ws.Cells.Item(7, 3).Value = 40

$ws.Cells.Item(8, 1).Value = "complex"
$ws.Cells.Item(8, 2).Value = 90
$ws.Cells.Item(8, 3).Value = 90

$ws.Cells.Item(9, 1).Value = "Superlist_website"
$ws.Cells.Item(9, 2).Value = 75
$ws.Cells.Item(9, 3).Value = 70

$ws.Cells.Item(10, 1).Value = "Spotify_website"
$ws.Cells.Item(10, 2).Value = 70
$ws.Cells.Item(10, 3).Value = 50

$ws.Cells.Item(11, 1).Value = "blog"
$ws.Cells.Item(11, 2).Value = 90
$ws.Cells.Item(11, 3).Value = 90

$ws.Cells.Item(12, 1).Value = "mubasic_website"
$ws.Cells.Item(12, 2).Value = 90
$ws.Cells.Item(12, 3).Value = 85

$ws.Cells.Item(13, 1).Value = "overflow_website"
$ws.Cells.Item(13, 2).Value = 80
$ws.Cells.Item(13, 3).Value = 75

$ws.Cells.Item(14, 1).Value = "RCA_website"
$ws.Cells.Item(14, 2).Value = 20
$ws.Cells.Item(14, 3).Value = 50

$ws.Cells.Item(15, 1).Value = "Crypto_website"
$ws.Cells.Item(15, 2).Value = 80
$ws.Cells.Item(15, 3).Value = 90

$ws.Cells.Item(16, 2).Value = 80
$ws.Cells.Item(16, 3).Value = 85
$ws.Cells.Item(16, 4).Value = "Crypto_website"

# Blank placeholder cells: exist in the sheet but carry no value, mirroring
# the authored file's empty inlineStr cells.
for ($r = 2; $r -le 15; $r++) {
    $blank.Copy($ws.Cells.Item($r, 4))
}
$blank.Copy($ws.Cells.Item(16, 1))

# D1 is blank too, but keeps the bold/bordered header style ("s=1") - copy it
# from A1 (same style) before A1 had a value, i.e. restore A1 after copying.
$ws.Cells.Item(1, 1).Copy($ws.Cells.Item(1, 4))
$ws.Cells.Item(1, 4).Value = ""
